$d = $word.ActiveDocument

# 1. Update experience years in the professional summary: "21 years" -> "15+ years"
$d.Content.Find.Execute("21 years of experience", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "15+ years of experience", 2) | Out-Null

# 2. Remove the EDUCATION section (heading + both degree entries) entirely.
#    Locate the "EDUCATION" heading paragraph and the paragraph containing the
#    second (last) degree entry, then delete the range spanning them.
$eduStart = $null
$eduEnd = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($text -match "^EDUCATION\s*$") {
        $eduStart = $i
    }
    elseif ($eduStart -ne $null -and $eduEnd -eq $null -and $text -match "Bachelor of Arts in Political Science") {
        $eduEnd = $i
    }
}

if ($eduStart -ne $null -and $eduEnd -ne $null) {
    $startRange = $d.Paragraphs.Item($eduStart).Range.Start
    $endRange = $d.Paragraphs.Item($eduEnd).Range.End
    $r = $d.Range($startRange, $endRange)
    $r.Delete()
}
